# Apply coin price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.576.41"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").Value = "1.850.49"
$ws.Range("E3").Value = "  -3.67%  "
$ws.Range("E4").Value = "  -1.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4663"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3921"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07880"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9851"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.20%  "
$ws.Range("D13").Value = "1.905.46"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.856"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.027"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06822"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001009"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").Value = "28.603.46"
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.416"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.59%  "
$ws.Range("E24").Value = "  -5.05%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.187.67"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.122"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.184"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.43%  "
$ws.Range("E29").Value = "  -3.33%  "
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9760"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09451"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.383"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.498"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.348"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06119"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02199"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.161"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5706"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.604"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1794"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.389"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.258"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5394"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07159"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.910"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.82%  "
